$d = $word.ActiveDocument

# 1. Update the letter date from September 19 to September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address paragraph "959 Story Road, San Jose CA 95122"
#    into two separate paragraphs: "959 Story Road" and "San Jose, CA 95122"
$addrRange = $d.Content
$addrRange.Find.Execute("959 Story Road, San Jose CA 95122", $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 0) | Out-Null

$addrPara = $addrRange.Paragraphs(1)
$addrParaRange = $addrPara.Range
$addrParaRange.Text = "959 Story Road`r"

$cityPara = $addrPara.Next()
$cityPara.Range.Text = "San Jose, CA 95122"

# 3. Remove the empty "NoSpacing" paragraph directly following "Board of Directors"
$bodRange = $d.Content
$bodRange.Find.Execute("Board of Directors", $false, $false, $false, $false, $false,
                        $true, 1, $false, "", 0) | Out-Null
$bodPara = $bodRange.Paragraphs(1)
$emptyPara = $bodPara.Next()
$emptyPara.Range.Delete() | Out-Null
